{"js": "// Locate the first table in the document body (the assignment info table).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\nconst table = tables.items[0];\n\n// Load the table's rows so we can reach the \"Student Name\" row (row index 1,\n// 0-based: row 0 = Assignment Date, row 1 = Student Name, row 2 = Roll Number).\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst studentNameRow = rows.items[1];\n\n// Set the row height to 197 twips (9.85 points) -> <w:trHeight w:val=\"197\"/>.\nstudentNameRow.preferredHeight = 9.85;\nawait context.sync();\n\n// Update the \"Student Name\" value cell (second cell) text: KIRTHIKA K -> PRAVEEN V\nconst nameRowCells = studentNameRow.cells;\nnameRowCells.load(\"items\");\nawait context.sync();\nconst nameCell = nameRowCells.items[1];\nconst nameRange = nameCell.body.getRange();\nconst nameSearchResults = nameRange.search(\"KIRTHIKA K\", { matchCase: true });\nnameSearchResults.load(\"items\");\nawait context.sync();\n\nif (nameSearchResults.items.length > 0) {\n  nameSearchResults.items[0].insertText(\"PRAVEEN V\", \"Replace\");\n  await context.sync();\n}\n\n// Update the \"Student Roll Number\" value cell. The roll number is split\n// across two runs: \"713319CS0\" + \"65\" -> \"713319CS\" + \"109\".\nconst rollRow = rows.items[2];\nconst rollRowCells = rollRow.cells;\nrollRowCells.load(\"items\");\nawait context.sync();\nconst rollCell = rollRowCells.items[1];\nconst rollRange = rollCell.body.getRange();\n\nconst rollPart1 = rollRange.search(\"713319CS0\", { matchCase: true });\nrollPart1.load(\"items\");\nawait context.sync();\nif (rollPart1.items.length > 0) {\n  rollPart1.items[0].insertText(\"713319CS\", \"Replace\");\n  await context.sync();\n}\n\nconst rollPart2Range = rollCell.body.getRange();\nconst rollPart2 = rollPart2Range.search(\"65\", { matchCase: true });\nrollPart2.load(\"items\");\nawait context.sync();\nif (rollPart2.items.length > 0) {\n  rollPart2.items[0].insertText(\"109\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Edit the assignment info table:\n#   - give the \"Student Name\" row an explicit height (197 twips / 9.85 pt)\n#   - change the student name from \"KIRTHIKA K\" to \"PRAVEEN V\"\n#   - change the roll number from \"713319CS0\" + \"65\" to \"713319CS\" + \"109\"\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Row 2 = \"Student Name\" row (Row 1 = Assignment Date, Row 2 = Student Name,\n# Row 3 = Student Roll Number, Row 4 = Maximum Marks).\n$nameRow = $tbl.Rows.Item(2)\n$nameRow.Height = 9.85\n\n# Student Name value cell (column 2 of row 2): KIRTHIKA K -> PRAVEEN V\n$nameCell = $tbl.Cell(2, 2)\n$nameCell.Range.Find.Execute(\"KIRTHIKA K\", $false, $false, $false, $false, $false, $true, 1, $false, \"PRAVEEN V\", 2)\n\n# Student Roll Number value cell (column 2 of row 3): \"713319CS0\"+\"65\" -> \"713319CS\"+\"109\"\n$rollCell = $tbl.Cell(3, 2)\n$rollCell.Range.Find.Execute(\"713319CS0\", $false, $false, $false, $false, $false, $true, 1, $false, \"713319CS\", 2)\n$rollCell.Range.Find.Execute(\"65\", $false, $false, $false, $false, $false, $true, 1, $false, \"109\", 2)\n"}
